{"js": "// Replace each three-digit-by-one-digit multiplication expression\n// in the document body with its updated version, matched by exact\n// (unique) original text so formatting/run properties are preserved.\nconst replacements = [\n  [\"364\u00d77=\", \"334\u00d76=\"],\n  [\"171\u00d73=\", \"291\u00d76=\"],\n  [\"748\u00d79=\", \"239\u00d73=\"],\n  [\"716\u00d77=\", \"945\u00d77=\"],\n  [\"318\u00d75=\", \"376\u00d77=\"],\n  [\"169\u00d72=\", \"831\u00d72=\"],\n  [\"675\u00d75=\", \"713\u00d76=\"],\n  [\"201\u00d77=\", \"239\u00d73=\"],\n  [\"568\u00d75=\", \"536\u00d76=\"],\n  [\"746\u00d72=\", \"464\u00d77=\"],\n  [\"520\u00d78=\", \"502\u00d77=\"],\n  [\"588\u00d78=\", \"401\u00d78=\"],\n  [\"699\u00d76=\", \"381\u00d77=\"],\n  [\"128\u00d77=\", \"182\u00d74=\"],\n  [\"792\u00d74=\", \"382\u00d74=\"],\n  [\"878\u00d72=\", \"307\u00d77=\"],\n  [\"860\u00d77=\", \"859\u00d78=\"],\n  [\"161\u00d74=\", \"715\u00d74=\"],\n  [\"880\u00d75=\", \"304\u00d79=\"],\n  [\"328\u00d76=\", \"379\u00d73=\"],\n  [\"289\u00d72=\", \"507\u00d75=\"],\n  [\"520\u00d72=\", \"660\u00d73=\"],\n  [\"484\u00d79=\", \"288\u00d72=\"],\n  [\"503\u00d75=\", \"299\u00d75=\"],\n  [\"983\u00d77=\", \"740\u00d72=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication expression\n# in the table to its new value. Matched via Find/Replace on the\n# exact original (unique) expression text so cell formatting is kept.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"364\u00d77=\"; New = \"334\u00d76=\" }\n    @{ Old = \"171\u00d73=\"; New = \"291\u00d76=\" }\n    @{ Old = \"748\u00d79=\"; New = \"239\u00d73=\" }\n    @{ Old = \"716\u00d77=\"; New = \"945\u00d77=\" }\n    @{ Old = \"318\u00d75=\"; New = \"376\u00d77=\" }\n    @{ Old = \"169\u00d72=\"; New = \"831\u00d72=\" }\n    @{ Old = \"675\u00d75=\"; New = \"713\u00d76=\" }\n    @{ Old = \"201\u00d77=\"; New = \"239\u00d73=\" }\n    @{ Old = \"568\u00d75=\"; New = \"536\u00d76=\" }\n    @{ Old = \"746\u00d72=\"; New = \"464\u00d77=\" }\n    @{ Old = \"520\u00d78=\"; New = \"502\u00d77=\" }\n    @{ Old = \"588\u00d78=\"; New = \"401\u00d78=\" }\n    @{ Old = \"699\u00d76=\"; New = \"381\u00d77=\" }\n    @{ Old = \"128\u00d77=\"; New = \"182\u00d74=\" }\n    @{ Old = \"792\u00d74=\"; New = \"382\u00d74=\" }\n    @{ Old = \"878\u00d72=\"; New = \"307\u00d77=\" }\n    @{ Old = \"860\u00d77=\"; New = \"859\u00d78=\" }\n    @{ Old = \"161\u00d74=\"; New = \"715\u00d74=\" }\n    @{ Old = \"880\u00d75=\"; New = \"304\u00d79=\" }\n    @{ Old = \"328\u00d76=\"; New = \"379\u00d73=\" }\n    @{ Old = \"289\u00d72=\"; New = \"507\u00d75=\" }\n    @{ Old = \"520\u00d72=\"; New = \"660\u00d73=\" }\n    @{ Old = \"484\u00d79=\"; New = \"288\u00d72=\" }\n    @{ Old = \"503\u00d75=\"; New = \"299\u00d75=\" }\n    @{ Old = \"983\u00d77=\"; New = \"740\u00d72=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
